$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
    2 = @(1.02,1.058160485012468,1.063598484966932,1.071366213069762,1.077708658268236,1.0552049193028,1.06315314894239,1.066316690672659,1.074063577813555,1.080389246855213,1.064662948359259)
    3 = @(1.02,1.059093648101666,1.064334487521868,1.072249595881637,1.078587871538276,1.055491230870675,1.063739145959732,1.066867933516909,1.074763334994814,1.081086054046755,1.06524977755958)
    4 = @(1.02,1.059698049085316,1.064811203255671,1.072822083231227,1.079157618058996,1.055675516836712,1.064118243528904,1.067224426321635,1.075216361554257,1.081537127217708,1.065629413490784)
    5 = @(1.02,1.059952277382248,1.065011726274564,1.073062965951771,1.079397338239582,1.055752756138058,1.064277595508309,1.067374247351922,1.075406869352822,1.081726803116696,1.065788991768277)
    6 = @(1.02,1.059994971496095,1.065045401480312,1.073103423433902,1.079437599921367,1.055765711184625,1.064304350165127,1.067399400083481,1.075438859664074,1.081758653120304,1.065815784419777)
    7 = @(1.02,1.059701445555055,1.064813882216922,1.072825301098479,1.079160820432831,1.055676549834381,1.064120372880303,1.067226428429701,1.075218906912569,1.081539661502198,1.065631545866106)
    8 = @(1.02,1.058475730589371,1.063847121676166,1.071664573374112,1.078005618308795,1.055301881473691,1.06335120563993,1.066503026144051,1.07430001421656,1.080624695501194,1.064861286320028)
    9 = @(1.02,1.056320376323778,1.062147253084663,1.069626027552205,1.075976482639574,1.054634217196285,1.061995256952316,1.065226829769313,1.07268267838227,1.079013939154782,1.063503412029747)
    10 = @(1.02,1.054886579342632,1.061016574120773,1.068271663375814,1.074628174406314,1.054184138472289,1.061090966820719,1.064375109653183,1.071605787612268,1.077941208028151,1.062597837702439)
    11 = @(1.02,1.054266481021629,1.060527604858273,1.067686334628885,1.074045416821102,1.053988079161572,1.060699335645377,1.064006100245446,1.07113981398687,1.077476982196264,1.062205650365908)
    12 = @(1.02,1.05403626202427,1.060346074851123,1.067469087243896,1.073829116840056,1.053915078532337,1.060553857289122,1.063869003399449,1.070966781040445,1.077304590473015,1.0620599654137)
    13 = @(1.02,1.054085639642601,1.060385009333643,1.06751567983835,1.073875506558275,1.053930745345175,1.06058506328526,1.063898412493157,1.071003894893787,1.077341567142649,1.062091215725932)
    14 = @(1.02,1.054247448732932,1.060512597585829,1.067668373420127,1.074027534080687,1.053982048482609,1.060687310538615,1.063994768389724,1.07112550999405,1.077462731369865,1.062193608182115)
    15 = @(1.02,1.054347159679515,1.060591221544525,1.067762475517623,1.074121224785767,1.054013634806476,1.060750307271639,1.064054132451984,1.071200447787763,1.077537390308643,1.062256694377725)
    16 = @(1.02,1.054927749026571,1.061049038634759,1.068310533441649,1.074666872759162,1.054197125640089,1.061116956757994,1.064399595281054,1.071636719747174,1.077972023046226,1.062623864548407)
    17 = @(1.02,1.055292137726671,1.061336382816695,1.06865461659685,1.075009430561148,1.054311911201977,1.061346928943351,1.064616239795748,1.071910470099412,1.078244730894343,1.062854163320644)
    18 = @(1.02,1.055504751316208,1.061504045750204,1.068855422429631,1.075209341748139,1.054378750529285,1.061481061236938,1.064742584640428,1.07207017553302,1.078403823102411,1.062988486097474)
    19 = @(1.02,1.055577259203309,1.061561224598764,1.068923910293127,1.075277523722665,1.054401521800982,1.061526795732762,1.06478566150248,1.072124636261069,1.07845807386182,1.063034285541528)
    20 = @(1.02,1.055253034829977,1.06130554729552,1.0686176885539,1.074972666684141,1.05429960749636,1.061322255778453,1.064592997974899,1.071881096011726,1.078215469189999,1.062829455117023)
    21 = @(1.02,1.05419979685651,1.060475023392759,1.067623404267255,1.073982761276298,1.053966945826006,1.060657201513695,1.063966394808184,1.071089695976712,1.07742705034358,1.062163456398925)
    22 = @(1.02,1.053538239450062,1.059953390121631,1.066999241285274,1.073361307284919,1.053756773327458,1.060239003372456,1.063572249627977,1.070592403614995,1.076931586189298,1.061744664368358)
    23 = @(1.02,1.053888881060621,1.060229865156718,1.067330028204033,1.07369066229222,1.053868285730096,1.060460702708785,1.063781209664651,1.070855999585207,1.077194217359397,1.061966678543174)
    24 = @(1.02,1.055270703524001,1.061319480355221,1.068634374413581,1.074989278379547,1.054305167360993,1.061333404540098,1.064603500015944,1.071894368800407,1.078228691213405,1.062840619711188)
    25 = @(1.02,1.056877044938302,1.062586263882378,1.070152224717399,1.076500285544853,1.054807703230323,1.062345864057152,1.065556924506107,1.073100569025553,1.079430168929346,1.063854517036887)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $row).Value = $values[$i]
    }
}

